# "some list manipulation is added"
# Turn the lone "1.line" cell into a small reading-log table:
#   day | type of book | number of page
#    1  | kk           | 5
#    2  | pr           | 10
#       | self         | 10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A2").Value = " day"
$ws.Range("B2").Value = "type of book"
$ws.Range("C2").Value = "number of page"

# Data rows
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "kk"
$ws.Range("C3").Value = 5

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "pr"
$ws.Range("C4").Value = 10

$ws.Range("B5").Value = "self"
$ws.Range("C5").Value = 10

# Widen the new text columns to fit their contents
$ws.Range("B2:C5").Columns.EntireColumn.AutoFit() | Out-Null

# Leave the selection where the last entry was typed
$ws.Range("C5").Select() | Out-Null
